$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: replace "Example 1" entry with the real log entry
$ws.Range("A4").Value = "Generate a Maze"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "05/01/2020"
$ws.Range("D4").Value = "Implementing Hunt and Kill(Drunken walk) algorithm to generate mazes "

# Row 5: replace "Example 2" entry with the real log entry
$ws.Range("A5").Value = "Making UI Functions"
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = "05/01/2020"
$ws.Range("D5").Value = "Made Basic UI Functions hit a snag at clearing the maze but its fixed"

# Row 6: fill in a new log entry
$ws.Range("A6").Value = "Made Temporary perfabs"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "05/01/2020"
$ws.Range("D6").Value = "Made Quick and dirty test Prefabs and also started to make some extra assets"

# Update selection to match the authored state
$ws.Range("D6").Select()

$wb.Application.Calculate()
